$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: MZ001 series -> MZ002 series (HR / Croatian dimension level)
$ws.Range("C2").Value = "MZ002"
$ws.Range("E2").Value = "hrv"
$ws.Range("F2").Value = "hr"
$ws.Range("J2").Value = "MZ002"
$ws.Range("K2").Value = "UMAR-SURS--MZ002--HR--M"

# Row 3: MZ001 series -> MZ002 series (SI / Slovenian dimension level)
$ws.Range("C3").Value = "MZ002"
$ws.Range("E3").Value = "slo"
$ws.Range("F3").Value = "si"
$ws.Range("J3").Value = "MZ002"
$ws.Range("K3").Value = "UMAR-SURS--MZ002--SI--M"

# Row 4: MZ002 series -> MZ007 series, interval switched from Monthly to Annual
$ws.Range("C4").Value = "MZ007"
$ws.Range("D4").Value = "sdfd--dsfg"
$ws.Range("E4").Value = "sdf--sdfSD"
$ws.Range("F4").Value = "LKJ--11"
$ws.Range("H4").Value = "A"
$ws.Range("I4").Value = "LKJ"
$ws.Range("J4").Value = "MZ007"
$ws.Range("K4").Value = "UMAR--MZ007--LKJ--11--A"

# The old fifth data row (former MZ007 row) is no longer present.
$ws.Rows.Item(5).Delete()

# Update the recorded selection to D5 (matches the saved workbook view).
$ws.Range("D5").Select()
